# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (default Office colours) - used by the notes master
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" colours - used by the slide master
#     (and therefore by every slide, since the slide master owns the active Design)
#
# The target edit swaps the two themes' contents: the slide master's theme
# becomes the plain "Office Theme" colour scheme, while the notes master's
# theme becomes the "Integral"/"Red Violet" colour scheme.
#
# The slide master's (active) theme colours are reachable through the
# PowerPoint object model via Slide.ThemeColorScheme, so re-point every
# theme colour slot to the "Office Theme" RGB values - this rewrites
# ppt/theme/theme2.xml (the part the slide master/presentation actually
# reference) to match the "Office Theme" colour scheme from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order is the standard MSO theme colour order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeRGB[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $comRGB = $b * 0x10000 + $g * 0x100 + $r
    $tcs.Item($i).RGB = $comRGB
}
